$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'schubert-winterreise_0'
$ws.Range("B2").Value = 'schubert-winterreise_36'
$ws.Range("C2").Value = 0.06153846153846154
$ws.Range("D2").Value = '[[''B:min'', ''E:min/B'', ''B:min''], [''B:min/F#'', ''F#:7'', ''B:min'']]'
$ws.Range("E2").Value = '[[''G:min/A#'', ''C:min'', ''G:min/A#''], [''G:min/D'', ''D:7'', ''G:min'']]'
$ws.Range("F2").Value = '[(25.48, 32.58), (74.1, 80.04)]'
$ws.Range("G2").Value = '[(97.88, 108.98), (47.58, 50.04)]'
$ws.Range("H2").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
$ws.Range("I2").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 3
$ws.Range("A3").Value = 'schubert-winterreise_33'
$ws.Range("B3").Value = 'isophonics_265'
$ws.Range("C3").Value = 0.3296703296703297
$ws.Range("D3").Value = '[[''G:maj/B'', ''C:maj'', ''G:maj/D'']]'
$ws.Range("E3").Value = '[[''A'', ''D'', ''A'']]'
$ws.Range("F3").Value = '[(63.6, 66.04)]'
$ws.Range("G3").Value = '[(105.425056, 115.003287)]'
$ws.Range("H3").Value = ''
$ws.Range("I3").Value = ''

# Row 4
$ws.Range("A4").Value = 'schubert-winterreise_170'
$ws.Range("B4").Value = 'schubert-winterreise_67'
$ws.Range("C4").Value = 0.2125874125874126
$ws.Range("D4").Value = '[[''G:min'', ''D:7'', ''G:min'', ''G:maj'']]'
$ws.Range("E4").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D'', ''B:maj/D#'']]'
$ws.Range("F4").Value = '[(48.8, 60.04)]'
$ws.Range("G4").Value = '[(0.3, 3.64)]'
$ws.Range("H4").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I4").Value = ''

# Row 5
$ws.Range("A5").Value = 'schubert-winterreise_210'
$ws.Range("B5").Value = 'schubert-winterreise_178'
$ws.Range("C5").Value = 0.3342175066312997
$ws.Range("D5").Value = '[[''D:maj/G'', ''G:min'', ''D:maj/G'', ''G:min'', ''D:maj/G'', ''G:min'']]'
$ws.Range("E5").Value = '[[''A:maj'', ''D:min'', ''A:maj'', ''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range("F5").Value = '[(36.4, 55.36)]'
$ws.Range("G5").Value = '[(1.14, 9.88)]'
$ws.Range("H5").Value = ''
$ws.Range("I5").Value = ''

# Row 6
$ws.Range("A6").Value = 'schubert-winterreise_170'
$ws.Range("B6").Value = 'schubert-winterreise_138'
$ws.Range("C6").Value = 0.1641025641025641
$ws.Range("D6").Value = '[[''C:min/D#'', ''G:min/D'', ''D:7'', ''G:min'']]'
$ws.Range("E6").Value = '[[''D:min'', ''A:min'', ''E:7'', ''A:min'']]'
$ws.Range("F6").Value = '[(40.46, 43.62)]'
$ws.Range("G6").Value = '[(16.72, 26.46)]'
$ws.Range("H6").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I6").Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'

# Row 7
$ws.Range("A7").Value = 'isophonics_159'
$ws.Range("B7").Value = 'isophonics_275'
$ws.Range("C7").Value = 0.5369318181818181
$ws.Range("D7").Value = '[[''A'', ''D/5'', ''A'', ''E/4'', ''D/5'', ''A'', ''D'']]'
$ws.Range("E7").Value = '[[''G'', ''C'', ''G'', ''D'', ''C'', ''G'', ''C'']]'
$ws.Range("F7").Value = '[(17.913, 30.076)]'
$ws.Range("G7").Value = '[(8.158789, 28.461467)]'
$ws.Range("H7").Value = ''
$ws.Range("I7").Value = ''

# Row 8
$ws.Range("A8").Value = 'schubert-winterreise_26'
$ws.Range("B8").Value = 'schubert-winterreise_147'
$ws.Range("C8").Value = 0.3939393939393939
$ws.Range("D8").Value = '[[''F:maj'', ''C:7'', ''F:maj'', ''C:7'', ''F:maj'']]'
$ws.Range("E8").Value = '[[''A:maj/E'', ''E:7'', ''A:maj'', ''E:7'', ''A:maj'']]'
$ws.Range("F8").Value = '[(62.5, 72.76)]'
$ws.Range("G8").Value = '[(19.78, 25.82)]'
$ws.Range("H8").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I8").Value = 'spotify:track:0XfunCHFEeQnzm4NaY8rJr'

# Row 9
$ws.Range("A9").Value = 'schubert-winterreise_9'
$ws.Range("B9").Value = 'jaah_87'
$ws.Range("C9").Value = 0.07964046822742475
$ws.Range("D9").Value = '[[''C:7'', ''F:min'', ''F:min'']]'
$ws.Range("E9").Value = '[[''C:7'', ''F:min'', ''F:min'']]'
$ws.Range("F9").Value = '[(5.68, 9.64)]'
$ws.Range("G9").Value = '[(55.46, 59.63)]'
$ws.Range("H9").Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'
$ws.Range("I9").Value = ''

# Row 10
$ws.Range("A10").Value = 'schubert-winterreise_203'
$ws.Range("B10").Value = 'jaah_67'
$ws.Range("C10").Value = 0.05274725274725275
$ws.Range("D10").Value = '[[''A:min7/C'', ''D:7'', ''G:maj'', ''D:7/C'']]'
$ws.Range("E10").Value = '[[''G:min7'', ''C:7'', ''F'', ''C:7'']]'
$ws.Range("F10").Value = '[(65.76, 75.0)]'
$ws.Range("G10").Value = '[(22.68, 27.0)]'
$ws.Range("H10").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I10").Value = ''

# Row 11
$ws.Range("A11").Value = 'jaah_86'
$ws.Range("B11").Value = 'jaah_77'
$ws.Range("C11").Value = 0.05595439189189189
$ws.Range("D11").Value = '[[''C:7'', ''C:7'', ''F:7'', ''F:7'', ''Bb'']]'
$ws.Range("E11").Value = '[[''D:7'', ''D:7'', ''G:7'', ''G:7'', ''C'']]'
$ws.Range("F11").Value = '[(26.8, 31.57)]'
$ws.Range("G11").Value = '[(7.47, 12.62)]'
$ws.Range("H11").Value = ''
$ws.Range("I11").Value = ''

# Row 12
$ws.Range("A12").Value = 'isophonics_111'
$ws.Range("B12").Value = 'jaah_5'
$ws.Range("C12").Value = 0.1391058667109855
$ws.Range("D12").Value = '[[''A:7'', ''E:7'', ''A:7''], [''A:7'', ''D:7'', ''A:7'']]'
$ws.Range("E12").Value = '[[''Bb:7'', ''F:7'', ''Bb:7''], [''Bb:7'', ''Eb:7'', ''Bb:7'']]'
$ws.Range("F12").Value = '[(11.993129, 22.326009), (1.834399, 15.441292)]'
$ws.Range("G12").Value = '[(22.23, 25.41), (11.55, 14.78)]'
$ws.Range("H12").Value = ''
$ws.Range("I12").Value = ''

# Row 13
$ws.Range("A13").Value = 'schubert-winterreise_70'
$ws.Range("B13").Value = 'schubert-winterreise_169'
$ws.Range("C13").Value = 0.2363636363636364
$ws.Range("D13").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D'']]'
$ws.Range("E13").Value = '[[''G:min'', ''D:7'', ''G:min'']]'
$ws.Range("F13").Value = '[(1.54, 3.58)]'
$ws.Range("G13").Value = '[(19.44, 28.3)]'
$ws.Range("H13").Value = ''
$ws.Range("I13").Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'

# Row 14
$ws.Range("A14").Value = 'isophonics_287'
$ws.Range("B14").Value = 'isophonics_112'
$ws.Range("C14").Value = 0.084375
$ws.Range("D14").Value = '[[''D'', ''A'', ''D''], [''A/3'', ''D'', ''A'']]'
$ws.Range("E14").Value = '[[''D/5'', ''A'', ''D/5''], [''A'', ''D/5'', ''A'']]'
$ws.Range("F14").Value = '[(21.722199, 26.760929), (65.81034, 71.679931)]'
$ws.Range("G14").Value = '[(1.922018, 6.206099), (0.421247, 4.824512)]'
$ws.Range("H14").Value = ''
$ws.Range("I14").Value = ''

# Row 15
$ws.Range("A15").Value = 'schubert-winterreise_40'
$ws.Range("B15").Value = 'isophonics_151'
$ws.Range("C15").Value = 0.2875
$ws.Range("D15").Value = '[[''D:maj/F#'', ''G:maj'', ''D:maj'']]'
$ws.Range("E15").Value = '[[''F'', ''Bb'', ''F'']]'
$ws.Range("F15").Value = '[(60.04, 67.08)]'
$ws.Range("G15").Value = '[(23.586235, 28.578526)]'
$ws.Range("H15").Value = ''
$ws.Range("I15").Value = ''

# Row 16
$ws.Range("A16").Value = 'schubert-winterreise_130'
$ws.Range("B16").Value = 'isophonics_79'
$ws.Range("C16").Value = 0.2965116279069767
$ws.Range("D16").Value = '[[''E:maj/G#'', ''A:maj'', ''E:maj'', ''B:maj'']]'
$ws.Range("E16").Value = '[[''E'', ''A'', ''E'', ''B'']]'
$ws.Range("F16").Value = '[(55.58, 66.2)]'
$ws.Range("G16").Value = '[(7.284457, 22.098788)]'
$ws.Range("H16").Value = ''
$ws.Range("I16").Value = ''

# Row 17
$ws.Range("A17").Value = 'isophonics_76'
$ws.Range("B17").Value = 'isophonics_22'
$ws.Range("C17").Value = 0.07023809523809524
$ws.Range("D17").Value = '[[''F:(1,5,9)'', ''C'', ''Bb'']]'
$ws.Range("E17").Value = '[[''F#:sus4'', ''F#'', ''E'']]'
$ws.Range("F17").Value = '[(118.453, 123.653)]'
$ws.Range("G17").Value = '[(31.005011, 35.323922)]'
$ws.Range("H17").Value = ''
$ws.Range("I17").Value = 'spotify:track:1h04XMpzGzmAudoI6VHBgA'
